$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) The document currently ends with two paragraphs:
#      - bold:   "Play Big Bad Wolf Free Slot Game - Review and Bonuses"
#      - italic: "Immerse yourself in the beautiful countryside ..."
#    We cut the bold paragraph and re-insert it right after the
#    Heading1 title, turning it into the new
#    "Meta description: Immerse yourself ..." paragraph.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldPara = $d.Paragraphs.Item($count - 1)

if ($boldPara.Range.Text -notmatch "Play Big Bad Wolf Free Slot Game - Review and Bonuses") {
    throw "Unexpected document layout: bold paragraph not found where expected."
}

$boldPara.Range.Cut()

$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"
$metaPara.Range.Paste()

# Turn the pasted "Play Big Bad Wolf Free Slot Game - Review and
# Bonuses" bold run into "Meta description" (still bold).
$metaPara.Range.Find.Execute("Play Big Bad Wolf Free Slot Game - Review and Bonuses", `
    $true, $false, $false, $false, $false, $true, 1, $false, "Meta description", 2)

# Append the (non-bold) description text right after the bold run,
# before the paragraph mark.
$insertPoint = $d.Range($metaPara.Range.End - 1, $metaPara.Range.End - 1)
$insertPoint.InsertAfter(": Immerse yourself in the beautiful countryside and win big with Big Bad Wolf slot game. Read our review and play for free now. Bonuses included.")

# ------------------------------------------------------------------
# 2) The old italic paragraph (now the last paragraph) keeps its
#    formatting and leading empty run, but its text is replaced with
#    the DALLE feature-image prompt.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

if ($lastPara.Range.Text -notmatch "Immerse yourself in the beautiful countryside") {
    throw "Unexpected document layout: italic paragraph not found where expected."
}

# Exclude the trailing paragraph mark so we only replace the run text
# (and so straight quotes/apostrophes aren't auto-corrected to curly ones).
$textRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$textRange.Text = "Prompt: Create a feature image fitting the game ""Big Bad Wolf"". DALLE, please create a cartoon-style feature image for the game ""Big Bad Wolf"" that features a happy Maya warrior with glasses. The Maya warrior should be holding a basket of apples and standing in front of a countryside landscape with hills and a straw house in the background. The image should also include the Wolf and the Three Little Pigs as cartoon characters. The setting should be under the moonlight, with stars shining brightly in the sky. The colors should be fun and vibrant, with a focus on shades of blue and yellow. Can't wait to see your creative work!"

Write-Host "Paragraphs:" $d.Paragraphs.Count
Write-Host "Meta paragraph:" $metaPara.Range.Text
Write-Host "Last paragraph:" $lastPara.Range.Text
